$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flat array of (row, productName, price) triples; nested arrays get
# flattened by this PowerShell engine, so we use a flat list and step by 3.
$data = @(
    2, 'Água Mineral sem Gás Frescca 510ml', 'R$ 1,09',
    3, 'Refrigerante Coca-Cola Sem Açúcar Pet 200ml', 'R$ 1,69',
    4, 'Água Mineral Cristal Select com Gás 500ml', 'R$ 1,69',
    5, 'Refrigerante Coca-Cola Pet 200ml', 'R$ 1,69',
    6, 'Refrigerante Guaraná Antarctica Zero Garrafa 200ML', 'R$ 1,65',
    7, 'Água Mineral sem Gás Frescca 1,5 Litros', 'R$ 2,29',
    8, 'Refrigerante Guaraná Antarctica 200ML Garrafa Pet', 'R$ 1,65',
    9, 'Refrigerante Coca-Cola Sem Açúcar Lata 220ml', 'R$ 2,69',
    10, 'Refrigerante Coca-Cola Sem Açúcar 2 Litros', 'R$ 10,29',
    11, 'Refrigerante Coca-Cola Sem Açúcar 1,5 litros', 'R$ 8,29',
    12, 'Suco de Maçã Yakult 200ml', 'R$ 2,97',
    13, 'Água Mineral com Gás Prata 510ml', 'R$ 2,39',
    14, 'Cerveja Budweiser American Lager 350ml Lata', 'R$ 3,78',
    15, 'Água Mineral Frescca com Gás 1,5L', 'R$ 2,39',
    16, 'Bebida Maguary Fruit Shoot 100% Suco Uva TP 150ml', 'R$ 2,18',
    17, 'Água Mineral Cristal Select sem Gás 500ml', 'R$ 1,29',
    18, 'Refrigerante Coca-Cola Lata 220ml', 'R$ 2,69',
    19, 'Refrigerante Coca Cola Original 2L', 'R$ 10,29',
    20, 'Refrigerante Guaraná Antarctica Sem Açúcar 350ml Lata', 'R$ 3,19',
    21, 'Refrigerante Coca-Cola Sem Açúcar 1 Litro', 'R$ 6,39',
    22, 'Água Mineral com Gás Crystal 500ml', 'R$ 2,59',
    23, 'Refrigerante Coca-Cola Sem Açúcar 600ml', 'R$ 4,79',
    24, 'Cerveja Petra Puro Malte 350ml', 'R$ 3,29',
    25, 'Cerveja Heineken Original Long Neck 250ml', 'R$ 5,29',
    26, 'Água Mineral Natural Prata Sem Gás 370Ml', 'R$ 2,39',
    27, 'Refrigerante Limoneto H2OH! 500ml', 'R$ 3,95',
    28, 'Suco de Morango Del Valle Kapo 200ml', 'R$ 2,39',
    29, 'Cerveja Heineken Premium Long Neck 330ml', 'R$ 5,98',
    30, 'Água Mineral Prata Com Gás 370Ml', 'R$ 2,89',
    31, 'Cerveja Amstel Puro Malte Lata 350ml', 'R$ 3,99',
    32, 'Água Mineral Sferriê com Gás 510ml', 'R$ 1,98',
    33, 'Água Mineral Natural Minalba sem Gás 510ml', 'R$ 1,75',
    34, 'Refrigerante Coca-Cola 600ml', 'R$ 4,79',
    35, 'Cerveja Budweiser Zero Álcool Lata 350ml', 'R$ 3,99',
    36, 'Refrigerante Pepsi Black Cola Zero 2 Litros', 'R$ 7,19',
    37, 'Suco Maguary Fruit Shoot 100% Maçã 150ml', 'R$ 2,18',
    38, 'Cerveja Heineken Original Lata 473ml', 'R$ 6,48',
    39, 'Cerveja Heineken Lata 350ml', 'R$ 5,69'
)

for ($i = 0; $i -lt $data.Count; $i += 3) {
    $r = $data[$i]
    $name = $data[$i + 1]
    $price = $data[$i + 2]
    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $price
}

# Row 40 no longer exists in the updated sheet; delete the whole row so the
# used range / dimension shrinks to A1:B39, matching the target layout.
$ws.Rows.Item(40).Delete()
